# Update "想去人数" (column F) counts across all sheets to match the
# regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 288
$ws.Range("F3").Value = 984
$ws.Range("F4").Value = 1230
$ws.Range("F5").Value = 1088
$ws.Range("F6").Value = 3220
$ws.Range("F9").Value = 1143
$ws.Range("F10").Value = 703
$ws.Range("F11").Value = 558
$ws.Range("F12").Value = 277
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 632
$ws.Range("F16").Value = 1412
$ws.Range("F17").Value = 1412
$ws.Range("F18").Value = 245
$ws.Range("F19").Value = 8
$ws.Range("F21").Value = 280
$ws.Range("F23").Value = 471
$ws.Range("F24").Value = 23820
$ws.Range("F25").Value = 23831
$ws.Range("F26").Value = 696
$ws.Range("F27").Value = 627
$ws.Range("F28").Value = 14547
$ws.Range("F29").Value = 14559
$ws.Range("F30").Value = 370
$ws.Range("F31").Value = 1
$ws.Range("F33").Value = 895
$ws.Range("F34").Value = 184
$ws.Range("F36").Value = 437
$ws.Range("F37").Value = 1119
$ws.Range("F38").Value = 5240
$ws.Range("F39").Value = 645
$ws.Range("F40").Value = 382
$ws.Range("F41").Value = 6
$ws.Range("F42").Value = 290
$ws.Range("F45").Value = 32

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 57
$ws.Range("F17").Value = 382
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 54
$ws.Range("F20").Value = 62
$ws.Range("F21").Value = 445
$ws.Range("F30").Value = 72
$ws.Range("F34").Value = 789
$ws.Range("F35").Value = 476
$ws.Range("F36").Value = 2
$ws.Range("F37").Value = 54
$ws.Range("F43").Value = 786
$ws.Range("F44").Value = 20
$ws.Range("F46").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 682
$ws.Range("F5").Value = 499
$ws.Range("F6").Value = 491

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 288
$ws.Range("F4").Value = 499
$ws.Range("F6").Value = 984
$ws.Range("F7").Value = 1230
$ws.Range("F9").Value = 1088
$ws.Range("F12").Value = 703
$ws.Range("F14").Value = 491
$ws.Range("F18").Value = 558
$ws.Range("F20").Value = 34
$ws.Range("F21").Value = 74
$ws.Range("F22").Value = 632
$ws.Range("F23").Value = 1412
$ws.Range("F24").Value = 1412
$ws.Range("F25").Value = 245
$ws.Range("F27").Value = 57
$ws.Range("F28").Value = 8
$ws.Range("F29").Value = 382
$ws.Range("F30").Value = 280
$ws.Range("F31").Value = 471
$ws.Range("F32").Value = 54
$ws.Range("F33").Value = 23847
$ws.Range("F34").Value = 62
$ws.Range("F35").Value = 696
$ws.Range("F36").Value = 14572
$ws.Range("F37").Value = 370
$ws.Range("F38").Value = 895
$ws.Range("F40").Value = 184
$ws.Range("F43").Value = 438
$ws.Range("F44").Value = 5240
$ws.Range("F45").Value = 72
$ws.Range("F46").Value = 645
$ws.Range("F47").Value = 476
$ws.Range("F48").Value = 54
$ws.Range("F49").Value = 54
$ws.Range("F50").Value = 290
$ws.Range("F54").Value = 32
